# Updates the cryptocurrency price/volume snapshot in cryptos.xlsx
# (scheduled GitHub Actions refresh - Mon Sep  4 16:59:26 UTC 2023).
#
# Rows 12/13 and 28/29 also swapped rank order (Polkadot now above
# WrappedEther; Cosmos now above EthereumClassic), so Coin name (B) and
# Link (C) are updated for those rows along with Price (D) and
# Volume(1h) (E).
#
# All edited cells hold text (prices/links/percent strings), never real
# numbers -- force Text format first so Excel's COM layer doesn't
# re-parse price strings like "1.001" / "4.400" / "0.5520" into doubles
# and silently drop the trailing zeros / dotted grouping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.861.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.630.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.57"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5068"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2576"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06330"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.44"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07764"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.247"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.634.93"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.854.13"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5498"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.70"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₅7647"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.880.93"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.03"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.400"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.867"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.021"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.914"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.37"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1246"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.64%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.755"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04887"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.245"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.186"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.541"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.368"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8954"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5520"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.539"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.119.01"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01553"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.54%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.592"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7974"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.15"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.765.61"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4444"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.65"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.584"
